# Update NATMI LR-pair edge table (Enho-Gpr19) with newly recomputed TPM-based
# statistics. The "ECs" sending cluster row-block is new (it previously had no
# qualifying cells), and the original "MuSCs" sending-cluster block below it
# has been recomputed against the updated TPM matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: ECs -> ECs --------------------------------------------------
$ws.Cells.Item(2,1).Value  = "ECs"
$ws.Cells.Item(2,2).Value  = "Enho"
$ws.Cells.Item(2,3).Value  = "Gpr19"
$ws.Cells.Item(2,4).Value  = "ECs"
$ws.Cells.Item(2,5).Value  = 1
$ws.Cells.Item(2,6).Value  = 0.3333333333333333
$ws.Cells.Item(2,7).Value  = 0.07501633333333334
$ws.Cells.Item(2,8).Value  = 0.225049
$ws.Cells.Item(2,9).Value  = 0.1517413757181704
$ws.Cells.Item(2,10).Value = 0.1517413757181704
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 1.093908333333333
$ws.Cells.Item(2,14).Value = 3.281725
$ws.Cells.Item(2,15).Value = 0.1455464035672057
$ws.Cells.Item(2,16).Value = 0.1455464035672057
$ws.Cells.Item(2,17).Value = 0.08206099216944444
$ws.Cells.Item(2,18).Value = 0.738548929525
$ws.Cells.Item(2,19).Value = 0.02208541150811983
$ws.Cells.Item(2,20).Value = 0.02208541150811982

# ---- Row 3: ECs -> FAPs --------------------------------------------------
$ws.Cells.Item(3,1).Value  = "ECs"
$ws.Cells.Item(3,2).Value  = "Enho"
$ws.Cells.Item(3,3).Value  = "Gpr19"
$ws.Cells.Item(3,4).Value  = "FAPs"
$ws.Cells.Item(3,5).Value  = 1
$ws.Cells.Item(3,6).Value  = 0.3333333333333333
$ws.Cells.Item(3,7).Value  = 0.07501633333333334
$ws.Cells.Item(3,8).Value  = 0.225049
$ws.Cells.Item(3,9).Value  = 0.1517413757181704
$ws.Cells.Item(3,10).Value = 0.1517413757181704
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 2.727648
$ws.Cells.Item(3,14).Value = 8.182944000000001
$ws.Cells.Item(3,15).Value = 0.3629183035726166
$ws.Cells.Item(3,16).Value = 0.3629183035726166
$ws.Cells.Item(3,17).Value = 0.204618151584
$ws.Cells.Item(3,18).Value = 1.841563364256
$ws.Cells.Item(3,19).Value = 0.05506972265741344
$ws.Cells.Item(3,20).Value = 0.05506972265741343

# ---- Row 4: ECs -> MuSCs --------------------------------------------------
$ws.Cells.Item(4,1).Value  = "ECs"
$ws.Cells.Item(4,2).Value  = "Enho"
$ws.Cells.Item(4,3).Value  = "Gpr19"
$ws.Cells.Item(4,4).Value  = "MuSCs"
$ws.Cells.Item(4,5).Value  = 1
$ws.Cells.Item(4,6).Value  = 0.3333333333333333
$ws.Cells.Item(4,7).Value  = 0.07501633333333334
$ws.Cells.Item(4,8).Value  = 0.225049
$ws.Cells.Item(4,9).Value  = 0.1517413757181704
$ws.Cells.Item(4,10).Value = 0.1517413757181704
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 2.690583
$ws.Cells.Item(4,14).Value = 8.071749000000001
$ws.Cells.Item(4,15).Value = 0.3579867409509296
$ws.Cells.Item(4,16).Value = 0.3579867409509296
$ws.Cells.Item(4,17).Value = 0.201837671189
$ws.Cells.Item(4,18).Value = 1.816539040701
$ws.Cells.Item(4,19).Value = 0.05432140056075836
$ws.Cells.Item(4,20).Value = 0.05432140056075835

# ---- Row 5: ECs -> Resolving-Mac -----------------------------------------
$ws.Cells.Item(5,1).Value  = "ECs"
$ws.Cells.Item(5,2).Value  = "Enho"
$ws.Cells.Item(5,3).Value  = "Gpr19"
$ws.Cells.Item(5,4).Value  = "Resolving-Mac"
$ws.Cells.Item(5,5).Value  = 1
$ws.Cells.Item(5,6).Value  = 0.3333333333333333
$ws.Cells.Item(5,7).Value  = 0.07501633333333334
$ws.Cells.Item(5,8).Value  = 0.225049
$ws.Cells.Item(5,9).Value  = 0.1517413757181704
$ws.Cells.Item(5,10).Value = 0.1517413757181704
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.003734
$ws.Cells.Item(5,14).Value = 3.011202
$ws.Cells.Item(5,15).Value = 0.1335485519092481
$ws.Cells.Item(5,16).Value = 0.1335485519092481
$ws.Cells.Item(5,17).Value = 0.075296444322
$ws.Cells.Item(5,18).Value = 0.677667998898
$ws.Cells.Item(5,19).Value = 0.02026484099187879
$ws.Cells.Item(5,20).Value = 0.02026484099187879

# ---- Row 6 (new): MuSCs -> ECs --------------------------------------------
$ws.Cells.Item(6,1).Value  = "MuSCs"
$ws.Cells.Item(6,2).Value  = "Enho"
$ws.Cells.Item(6,3).Value  = "Gpr19"
$ws.Cells.Item(6,4).Value  = "ECs"
$ws.Cells.Item(6,5).Value  = 1
$ws.Cells.Item(6,6).Value  = 0.3333333333333333
$ws.Cells.Item(6,7).Value  = 0.4193533333333333
$ws.Cells.Item(6,8).Value  = 1.25806
$ws.Cells.Item(6,9).Value  = 0.8482586242818295
$ws.Cells.Item(6,10).Value = 0.8482586242818295
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 1.093908333333333
$ws.Cells.Item(6,14).Value = 3.281725
$ws.Cells.Item(6,15).Value = 0.1455464035672057
$ws.Cells.Item(6,16).Value = 0.1455464035672057
$ws.Cells.Item(6,17).Value = 0.4587341059444444
$ws.Cells.Item(6,18).Value = 4.128606953499999
$ws.Cells.Item(6,19).Value = 0.1234609920590859
$ws.Cells.Item(6,20).Value = 0.1234609920590859

# ---- Row 7 (new): MuSCs -> FAPs -------------------------------------------
$ws.Cells.Item(7,1).Value  = "MuSCs"
$ws.Cells.Item(7,2).Value  = "Enho"
$ws.Cells.Item(7,3).Value  = "Gpr19"
$ws.Cells.Item(7,4).Value  = "FAPs"
$ws.Cells.Item(7,5).Value  = 1
$ws.Cells.Item(7,6).Value  = 0.3333333333333333
$ws.Cells.Item(7,7).Value  = 0.4193533333333333
$ws.Cells.Item(7,8).Value  = 1.25806
$ws.Cells.Item(7,9).Value  = 0.8482586242818295
$ws.Cells.Item(7,10).Value = 0.8482586242818295
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 2.727648
$ws.Cells.Item(7,14).Value = 8.182944000000001
$ws.Cells.Item(7,15).Value = 0.3629183035726166
$ws.Cells.Item(7,16).Value = 0.3629183035726166
$ws.Cells.Item(7,17).Value = 1.14384828096
$ws.Cells.Item(7,18).Value = 10.29463452864
$ws.Cells.Item(7,19).Value = 0.3078485809152031
$ws.Cells.Item(7,20).Value = 0.3078485809152031

# ---- Row 8 (new): MuSCs -> MuSCs ------------------------------------------
$ws.Cells.Item(8,1).Value  = "MuSCs"
$ws.Cells.Item(8,2).Value  = "Enho"
$ws.Cells.Item(8,3).Value  = "Gpr19"
$ws.Cells.Item(8,4).Value  = "MuSCs"
$ws.Cells.Item(8,5).Value  = 1
$ws.Cells.Item(8,6).Value  = 0.3333333333333333
$ws.Cells.Item(8,7).Value  = 0.4193533333333333
$ws.Cells.Item(8,8).Value  = 1.25806
$ws.Cells.Item(8,9).Value  = 0.8482586242818295
$ws.Cells.Item(8,10).Value = 0.8482586242818295
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 2.690583
$ws.Cells.Item(8,14).Value = 8.071749000000001
$ws.Cells.Item(8,15).Value = 0.3579867409509296
$ws.Cells.Item(8,16).Value = 0.3579867409509296
$ws.Cells.Item(8,17).Value = 1.12830494966
$ws.Cells.Item(8,18).Value = 10.15474454694
$ws.Cells.Item(8,19).Value = 0.3036653403901712
$ws.Cells.Item(8,20).Value = 0.3036653403901712

# ---- Row 9 (new): MuSCs -> Resolving-Mac -----------------------------------
$ws.Cells.Item(9,1).Value  = "MuSCs"
$ws.Cells.Item(9,2).Value  = "Enho"
$ws.Cells.Item(9,3).Value  = "Gpr19"
$ws.Cells.Item(9,4).Value  = "Resolving-Mac"
$ws.Cells.Item(9,5).Value  = 1
$ws.Cells.Item(9,6).Value  = 0.3333333333333333
$ws.Cells.Item(9,7).Value  = 0.4193533333333333
$ws.Cells.Item(9,8).Value  = 1.25806
$ws.Cells.Item(9,9).Value  = 0.8482586242818295
$ws.Cells.Item(9,10).Value = 0.8482586242818295
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 1.003734
$ws.Cells.Item(9,14).Value = 3.011202
$ws.Cells.Item(9,15).Value = 0.1335485519092481
$ws.Cells.Item(9,16).Value = 0.1335485519092481
$ws.Cells.Item(9,17).Value = 0.4209191986799999
$ws.Cells.Item(9,18).Value = 3.78827278812
$ws.Cells.Item(9,19).Value = 0.1132837109173692
$ws.Cells.Item(9,20).Value = 0.1132837109173692

Write-Output "Updated Enho-Gpr19 LR-pair table with new TPM-derived statistics (rows 2-9)."
